$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.01532382679555688
$ws.Range("E2").Value = 0.002755130728141708
$ws.Range("D3").Value = 0.05029184251188324
$ws.Range("E3").Value = 0.003023534711954934
$ws.Range("D4").Value = 0.01451279157442988
$ws.Range("E4").Value = 0.004705571924886964
$ws.Range("D5").Value = 0.009943677942284953
$ws.Range("E5").Value = 0.0177150192554556
$ws.Range("D6").Value = 0.0156489110291835
$ws.Range("E6").Value = 0.01483741246619985
$ws.Range("D7").Value = 0.02032894964376547
$ws.Range("E7").Value = 0.01608381289656213
$ws.Range("D8").Value = 0.00473837315655404
$ws.Range("E8").Value = -0.009004337167641041
$ws.Range("D9").Value = 0.006897517089721562
$ws.Range("E9").Value = -0.01232511658894087
$ws.Range("D10").Value = 0.01437088035552955
$ws.Range("E10").Value = 0.02156404260846956
$ws.Range("D11").Value = 0.008487229093485578
$ws.Range("E11").Value = 0.00579034441269366
$ws.Range("D12").Value = 0.01567756771336891
$ws.Range("E12").Value = 0.02540302882266743
$ws.Range("D13").Value = 0.002767699576169365
$ws.Range("E13").Value = 0.01372074253430178
$ws.Range("D14").Value = 0.005795479303910046
$ws.Range("E14").Value = 0.0115081768625076
$ws.Range("D15").Value = 0.01457103038137683
$ws.Range("E15").Value = 0.02572245157192765
$ws.Range("D16").Value = 0.01065961637179446
$ws.Range("E16").Value = 0.01427994072477423
$ws.Range("D17").Value = 0.02101766741446635
$ws.Range("E17").Value = -0.001571467614102429
$ws.Range("D18").Value = 0.008648893639235307
$ws.Range("E18").Value = 0.007903300790330192
$ws.Range("D19").Value = 0.01717613997322032
$ws.Range("E19").Value = 0.01046373365041631
$ws.Range("D20").Value = 0.01233067123523397
$ws.Range("E20").Value = 0.009006211180124346
$ws.Range("D21").Value = 0.007549600257477179
$ws.Range("E21").Value = 0.007507037847982678
$ws.Range("D22").Value = 0.01490045460280882
$ws.Range("E22").Value = 0.0106740225473736
$ws.Range("D23").Value = 0.01997137932494871
$ws.Range("E23").Value = 0.0265757798006816
$ws.Range("D24").Value = 0.01018537696903356
$ws.Range("E24").Value = 0.004361257495911319
$ws.Range("D25").Value = 0.0199776021015814
$ws.Range("E25").Value = 0.01578204169781539
$ws.Range("D26").Value = 0.01411995885469455
$ws.Range("E26").Value = 0.02358803986710956
$ws.Range("D27").Value = 0.01957580220115999
$ws.Range("E27").Value = 0.04852332024315476
$ws.Range("D28").Value = 0.05484914882426357
$ws.Range("E28").Value = 0.01791968722000492
$ws.Range("D29").Value = 0.02064797065913461
$ws.Range("E29").Value = 0.01586931155192528
$ws.Range("D30").Value = 0.02854563155626681
$ws.Range("E30").Value = 0.01135804677368846
$ws.Range("D31").Value = 0.01471428189062881
$ws.Range("E31").Value = 0.01998074145402007
$ws.Range("D32").Value = 0.01331361464979835
$ws.Range("E32").Value = 0.006375808361417157
$ws.Range("D33").Value = 0.01764587982980066
$ws.Range("E33").Value = 0.005078125000000044
$ws.Range("D34").Value = 0.04212819780329743
$ws.Range("E34").Value = 0.01308487671855474
$ws.Range("D35").Value = 0.01095208687353079
$ws.Range("E35").Value = 0.01538461538461533
$ws.Range("D36").Value = 0.01027198525509056
$ws.Range("E36").Value = 0.01426583159359773
$ws.Range("D37").Value = 0.01029869532709856
$ws.Range("E37").Value = 0.02184522426214275
$ws.Range("D38").Value = 0.007578384588362742
$ws.Range("E38").Value = 0.01768569984840851
$ws.Range("D39").Value = 0.0122213418365489
$ws.Range("E39").Value = 0.01602719766877136
$ws.Range("D40").Value = 0.01745440977956412
$ws.Range("E40").Value = 0.017248062015504
$ws.Range("D41").Value = 0.01742447662837713
$ws.Range("E41").Value = 0.01468439000626343
$ws.Range("D42").Value = 0.03194390583121522
$ws.Range("E42").Value = -0.000239757844577082
$ws.Range("D43").Value = 0.01142905727744157
$ws.Range("E43").Value = 0.01405318459274407
$ws.Range("D44").Value = 0.02182608587824011
$ws.Range("E44").Value = 0.01409599782441195
$ws.Range("D45").Value = 0.01213096797283725
$ws.Range("E45").Value = 0.01946114872863292
$ws.Range("D46").Value = 0.00850155743557828
$ws.Range("E46").Value = 0.02496537279145961
$ws.Range("D47").Value = 0.01363190169664158
$ws.Range("E47").Value = 0.01158772964773291
$ws.Range("D48").Value = 0.01077708324761458
$ws.Range("E48").Value = 0.02006715701950768
$ws.Range("D49").Value = 0.01576755863698009
$ws.Range("E49").Value = 0.01880995749848213
$ws.Range("D50").Value = 0.008489016147287785
$ws.Range("E50").Value = 0.03610685071574649
$ws.Range("D51").Value = 0.01119876412158555
$ws.Range("E51").Value = 0.02837033026529512
$ws.Range("D52").Value = 0.008482250872179428
$ws.Range("E52").Value = 0.01009766594934614
$ws.Range("D53").Value = 0.01013594578439749
$ws.Range("E53").Value = 0.00484848484848488
$ws.Range("D54").Value = 0.1385214331309368
$ws.Range("E54").Value = 0.0003943995267203793
$ws.Range("D55").Value = 0.04364942735242684
$ws.Range("E55").Value = 0.01255282127765356
$ws.Range("E56").Value = 0.01199710168721446

$ws.Protect()
